$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtTest10mm")

# New "Pressure" data row (row 17), values observed during the knee ICR test
$ws.Range("R22").Value = "kmax"
$ws.Range("S22").Formula = "=1-440/520"

$ws.Range("R23").Value = "avg. pres."
$ws.Range("S23").Formula = "=SUM(C17:K17)/9"

$ws.Range("B17").Value = "Pressure"
$ws.Range("C17").Value = 601.48289999999997
$ws.Range("D17").Value = 600.72787878787904
$ws.Range("E17").Value = 585.37148582600196
$ws.Range("F17").Value = 597.24558162267795
$ws.Range("G17").Value = 595.76131964809395
$ws.Range("H17").Value = 597.24558162267795
$ws.Range("I17").Value = 597.24558162267795
$ws.Range("J17").Value = 596.503450635386
$ws.Range("K17").Value = 595.76131964809395

# Switch active sheet / selection to match the author's final cursor position
$null = $ws.Activate()
$null = $ws.Range("F19").Select()
